$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 1.726994276046753
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = 6606.738577097564
$ws.Cells.Item(2, 6).Value = 0.2127634481962593
$ws.Cells.Item(2, 7).Value = 0.1934638974912634
$ws.Cells.Item(2, 8).Value = 0.1733836405884461
$ws.Cells.Item(2, 9).Value = 0.1623690174035624
$ws.Cells.Item(2, 10).Value = 0.1469256233805147
$ws.Cells.Item(2, 11).Value = 0.1450049636073846
$ws.Cells.Item(2, 12).Value = 0.1450049636073846
$ws.Cells.Item(2, 13).Value = 0.1450049636073846
$ws.Cells.Item(2, 14).Value = 0.1448036855417751
$ws.Cells.Item(2, 15).Value = 0.1448036855417751
$ws.Cells.Item(2, 16).Value = 0.1448036855417751
$ws.Cells.Item(2, 17).Value = 0.1448036855417751
$ws.Cells.Item(2, 18).Value = 0.1448036855417751
$ws.Cells.Item(2, 19).Value = 0.1448036855417751
$ws.Cells.Item(2, 20).Value = 0.1447893721435916
$ws.Cells.Item(2, 21).Value = 0.1447893721435916
$ws.Cells.Item(2, 22).Value = 0.1447863270389388
$ws.Cells.Item(2, 23).Value = 0.1447863270389388
$ws.Cells.Item(2, 24).Value = 0.1447863270389388
$ws.Cells.Item(2, 25).Value = 0.1447863270389388

$ws.Cells.Item(3, 3).Value = 1.967994451522827
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(3, 5).Value = 6607.667621530278
$ws.Cells.Item(3, 6).Value = 0.2127634481962593
$ws.Cells.Item(3, 7).Value = 0.1934638974912634
$ws.Cells.Item(3, 8).Value = 0.1733836405884461
$ws.Cells.Item(3, 9).Value = 0.1622856921731826
$ws.Cells.Item(3, 10).Value = 0.1469737490904527
$ws.Cells.Item(3, 11).Value = 0.1448259699691833
$ws.Cells.Item(3, 12).Value = 0.1448259699691833
$ws.Cells.Item(3, 13).Value = 0.1448259699691833
$ws.Cells.Item(3, 14).Value = 0.1448259699691833
$ws.Cells.Item(3, 15).Value = 0.1448259699691833
$ws.Cells.Item(3, 16).Value = 0.1448259699691833
$ws.Cells.Item(3, 17).Value = 0.1448259699691833
$ws.Cells.Item(3, 18).Value = 0.1448259699691833
$ws.Cells.Item(3, 19).Value = 0.1448259699691833
$ws.Cells.Item(3, 20).Value = 0.1448259699691833
$ws.Cells.Item(3, 21).Value = 0.1448259699691833
$ws.Cells.Item(3, 22).Value = 0.144804437066867
$ws.Cells.Item(3, 23).Value = 0.144804437066867
$ws.Cells.Item(3, 24).Value = 0.144804437066867
$ws.Cells.Item(3, 25).Value = 0.144804437066867

$ws.Cells.Item(4, 3).Value = 1.651999711990356
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(4, 5).Value = 6606.94576317169
$ws.Cells.Item(4, 6).Value = 0.2127634481962593
$ws.Cells.Item(4, 7).Value = 0.1934638974912634
$ws.Cells.Item(4, 8).Value = 0.1733836405884461
$ws.Cells.Item(4, 9).Value = 0.1585552995096235
$ws.Cells.Item(4, 10).Value = 0.1461568894446195
$ws.Cells.Item(4, 11).Value = 0.1449084267808127
$ws.Cells.Item(4, 12).Value = 0.1449084267808127
$ws.Cells.Item(4, 13).Value = 0.1449084267808127
$ws.Cells.Item(4, 14).Value = 0.1449084267808127
$ws.Cells.Item(4, 15).Value = 0.1449084267808127
$ws.Cells.Item(4, 16).Value = 0.1449084267808127
$ws.Cells.Item(4, 17).Value = 0.1449084267808127
$ws.Cells.Item(4, 18).Value = 0.1449084267808127
$ws.Cells.Item(4, 19).Value = 0.1448900412810464
$ws.Cells.Item(4, 20).Value = 0.1448900412810464
$ws.Cells.Item(4, 21).Value = 0.1448786037712648
$ws.Cells.Item(4, 22).Value = 0.1448786037712648
$ws.Cells.Item(4, 23).Value = 0.1448621142022455
$ws.Cells.Item(4, 24).Value = 0.144808520661192
$ws.Cells.Item(4, 25).Value = 0.1447903657538341

$ws.Cells.Item(5, 3).Value = 1.791008472442627
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 6606.772291975225
$ws.Cells.Item(5, 6).Value = 0.2127634481962593
$ws.Cells.Item(5, 7).Value = 0.1934638974912634
$ws.Cells.Item(5, 8).Value = 0.1733836405884461
$ws.Cells.Item(5, 9).Value = 0.1616210852956254
$ws.Cells.Item(5, 10).Value = 0.1450589657350946
$ws.Cells.Item(5, 11).Value = 0.1447921023300285
$ws.Cells.Item(5, 12).Value = 0.1447921023300285
$ws.Cells.Item(5, 13).Value = 0.1447921023300285
$ws.Cells.Item(5, 14).Value = 0.1447921023300285
$ws.Cells.Item(5, 15).Value = 0.1447921023300285
$ws.Cells.Item(5, 16).Value = 0.1447921023300285
$ws.Cells.Item(5, 17).Value = 0.1447869842490297
$ws.Cells.Item(5, 18).Value = 0.1447869842490297
$ws.Cells.Item(5, 19).Value = 0.1447869842490297
$ws.Cells.Item(5, 20).Value = 0.1447869842490297
$ws.Cells.Item(5, 21).Value = 0.1447869842490297
$ws.Cells.Item(5, 22).Value = 0.1447869842490297
$ws.Cells.Item(5, 23).Value = 0.1447869842490297
$ws.Cells.Item(5, 24).Value = 0.1447869842490297
$ws.Cells.Item(5, 25).Value = 0.1447869842490297

$ws.Cells.Item(6, 3).Value = 1.742999315261841
$ws.Cells.Item(6, 4).Value = 2
$ws.Cells.Item(6, 5).Value = 6606.899230796525
$ws.Cells.Item(6, 6).Value = 0.2127634481962593
$ws.Cells.Item(6, 7).Value = 0.1934638974912634
$ws.Cells.Item(6, 8).Value = 0.1733836405884461
$ws.Cells.Item(6, 9).Value = 0.1623690174035624
$ws.Cells.Item(6, 10).Value = 0.1498334151688435
$ws.Cells.Item(6, 11).Value = 0.1455606811692347
$ws.Cells.Item(6, 12).Value = 0.1447894586899907
$ws.Cells.Item(6, 13).Value = 0.1447894586899907
$ws.Cells.Item(6, 14).Value = 0.1447894586899907
$ws.Cells.Item(6, 15).Value = 0.1447894586899907
$ws.Cells.Item(6, 16).Value = 0.1447894586899907
$ws.Cells.Item(6, 17).Value = 0.1447894586899907
$ws.Cells.Item(6, 18).Value = 0.1447894586899907
$ws.Cells.Item(6, 19).Value = 0.1447894586899907
$ws.Cells.Item(6, 20).Value = 0.1447894586899907
$ws.Cells.Item(6, 21).Value = 0.1447894586899907
$ws.Cells.Item(6, 22).Value = 0.1447894586899907
$ws.Cells.Item(6, 23).Value = 0.1447894586899907
$ws.Cells.Item(6, 24).Value = 0.1447894586899907
$ws.Cells.Item(6, 25).Value = 0.1447894586899907

$ws.Cells.Item(7, 3).Value = 1.979997158050537
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = 6607.456545467006
$ws.Cells.Item(7, 6).Value = 0.2127634481962593
$ws.Cells.Item(7, 7).Value = 0.1934638974912634
$ws.Cells.Item(7, 8).Value = 0.1733836405884461
$ws.Cells.Item(7, 9).Value = 0.1623690174035624
$ws.Cells.Item(7, 10).Value = 0.14497875291081
$ws.Cells.Item(7, 11).Value = 0.14497875291081
$ws.Cells.Item(7, 12).Value = 0.14497875291081
$ws.Cells.Item(7, 13).Value = 0.144847932583505
$ws.Cells.Item(7, 14).Value = 0.144847932583505
$ws.Cells.Item(7, 15).Value = 0.1448380603416486
$ws.Cells.Item(7, 16).Value = 0.1448380603416486
$ws.Cells.Item(7, 17).Value = 0.1448380603416486
$ws.Cells.Item(7, 18).Value = 0.1448190502397089
$ws.Cells.Item(7, 19).Value = 0.1448190502397089
$ws.Cells.Item(7, 20).Value = 0.1448190502397089
$ws.Cells.Item(7, 21).Value = 0.1448190502397089
$ws.Cells.Item(7, 22).Value = 0.1448190502397089
$ws.Cells.Item(7, 23).Value = 0.1448190502397089
$ws.Cells.Item(7, 24).Value = 0.1448190502397089
$ws.Cells.Item(7, 25).Value = 0.1448003225237233

$ws.Cells.Item(8, 3).Value = 1.858994960784912
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(8, 5).Value = 6606.74368701515
$ws.Cells.Item(8, 6).Value = 0.2127634481962593
$ws.Cells.Item(8, 7).Value = 0.1934638974912634
$ws.Cells.Item(8, 8).Value = 0.1733836405884461
$ws.Cells.Item(8, 9).Value = 0.1623690174035624
$ws.Cells.Item(8, 10).Value = 0.1465262001453601
$ws.Cells.Item(8, 11).Value = 0.145229681148064
$ws.Cells.Item(8, 12).Value = 0.145229681148064
$ws.Cells.Item(8, 13).Value = 0.145229681148064
$ws.Cells.Item(8, 14).Value = 0.1449128490878133
$ws.Cells.Item(8, 15).Value = 0.1448258264878176
$ws.Cells.Item(8, 16).Value = 0.1448258264878176
$ws.Cells.Item(8, 17).Value = 0.1448258264878176
$ws.Cells.Item(8, 18).Value = 0.1448258264878176
$ws.Cells.Item(8, 19).Value = 0.1447864266474688
$ws.Cells.Item(8, 20).Value = 0.1447864266474688
$ws.Cells.Item(8, 21).Value = 0.1447864266474688
$ws.Cells.Item(8, 22).Value = 0.1447864266474688
$ws.Cells.Item(8, 23).Value = 0.1447864266474688
$ws.Cells.Item(8, 24).Value = 0.1447864266474688
$ws.Cells.Item(8, 25).Value = 0.1447864266474688

$ws.Cells.Item(9, 3).Value = 1.858001232147217
$ws.Cells.Item(9, 4).Value = 2
$ws.Cells.Item(9, 5).Value = 6606.803488594145
$ws.Cells.Item(9, 6).Value = 0.2127634481962593
$ws.Cells.Item(9, 7).Value = 0.1934638974912634
$ws.Cells.Item(9, 8).Value = 0.1733836405884461
$ws.Cells.Item(9, 9).Value = 0.1620588333831717
$ws.Cells.Item(9, 10).Value = 0.1471539925302741
$ws.Cells.Item(9, 11).Value = 0.1449685875807807
$ws.Cells.Item(9, 12).Value = 0.1449685875807807
$ws.Cells.Item(9, 13).Value = 0.1449685875807807
$ws.Cells.Item(9, 14).Value = 0.1447875923702562
$ws.Cells.Item(9, 15).Value = 0.1447875923702562
$ws.Cells.Item(9, 16).Value = 0.1447875923702562
$ws.Cells.Item(9, 17).Value = 0.1447875923702562
$ws.Cells.Item(9, 18).Value = 0.1447875923702562
$ws.Cells.Item(9, 19).Value = 0.1447875923702562
$ws.Cells.Item(9, 20).Value = 0.1447875923702562
$ws.Cells.Item(9, 21).Value = 0.1447875923702562
$ws.Cells.Item(9, 22).Value = 0.1447875923702562
$ws.Cells.Item(9, 23).Value = 0.1447875923702562
$ws.Cells.Item(9, 24).Value = 0.1447875923702562
$ws.Cells.Item(9, 25).Value = 0.1447875923702562

$ws.Cells.Item(10, 3).Value = 2.090995788574219
$ws.Cells.Item(10, 4).Value = 2
$ws.Cells.Item(10, 5).Value = 6606.826419322614
$ws.Cells.Item(10, 6).Value = 0.2127634481962593
$ws.Cells.Item(10, 7).Value = 0.1934638974912634
$ws.Cells.Item(10, 8).Value = 0.1733836405884461
$ws.Cells.Item(10, 9).Value = 0.1603734124806439
$ws.Cells.Item(10, 10).Value = 0.1457985613705054
$ws.Cells.Item(10, 11).Value = 0.1451580564148081
$ws.Cells.Item(10, 12).Value = 0.1451224548961737
$ws.Cells.Item(10, 13).Value = 0.1448007549883789
$ws.Cells.Item(10, 14).Value = 0.1448007549883789
$ws.Cells.Item(10, 15).Value = 0.1448007549883789
$ws.Cells.Item(10, 16).Value = 0.1448007549883789
$ws.Cells.Item(10, 17).Value = 0.1448007549883789
$ws.Cells.Item(10, 18).Value = 0.1448007549883789
$ws.Cells.Item(10, 19).Value = 0.1448007549883789
$ws.Cells.Item(10, 20).Value = 0.1447886511353449
$ws.Cells.Item(10, 21).Value = 0.1447886511353449
$ws.Cells.Item(10, 22).Value = 0.1447886511353449
$ws.Cells.Item(10, 23).Value = 0.1447886511353449
$ws.Cells.Item(10, 24).Value = 0.1447880393630139
$ws.Cells.Item(10, 25).Value = 0.1447880393630139

$ws.Cells.Item(11, 3).Value = 1.782000780105591
$ws.Cells.Item(11, 4).Value = 2
$ws.Cells.Item(11, 5).Value = 6606.722414897682
$ws.Cells.Item(11, 6).Value = 0.2127634481962593
$ws.Cells.Item(11, 7).Value = 0.1934638974912634
$ws.Cells.Item(11, 8).Value = 0.1733836405884461
$ws.Cells.Item(11, 9).Value = 0.1555526384135345
$ws.Cells.Item(11, 10).Value = 0.145476336841981
$ws.Cells.Item(11, 11).Value = 0.145476336841981
$ws.Cells.Item(11, 12).Value = 0.1450182489554601
$ws.Cells.Item(11, 13).Value = 0.1450182489554601
$ws.Cells.Item(11, 14).Value = 0.1450182489554601
$ws.Cells.Item(11, 15).Value = 0.1448397747604809
$ws.Cells.Item(11, 16).Value = 0.1448006014206911
$ws.Cells.Item(11, 17).Value = 0.1448006014206911
$ws.Cells.Item(11, 18).Value = 0.1448006014206911
$ws.Cells.Item(11, 19).Value = 0.1448006014206911
$ws.Cells.Item(11, 20).Value = 0.1448006014206911
$ws.Cells.Item(11, 21).Value = 0.1447969382309313
$ws.Cells.Item(11, 22).Value = 0.1447969382309313
$ws.Cells.Item(11, 23).Value = 0.1447969382309313
$ws.Cells.Item(11, 24).Value = 0.1447969382309313
$ws.Cells.Item(11, 25).Value = 0.1447969382309313
